$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 3103.5
$ws.Range("J17").Value = 3103.5
$ws.Range("L17").Value = 9310.5
$ws.Range("N17").Value = -9646.5

$ws.Range("H40").Value = 1849.3
$ws.Range("I40").Value = 1465
$ws.Range("K40").Value = 1465
$ws.Range("M40").Value = -1290

$ws.Range("H64").Value = 4292
$ws.Range("I64").Value = 4990
$ws.Range("J64").Value = 3826.6667
$ws.Range("K64").Value = 4990
$ws.Range("L64").Value = 3826.6667
$ws.Range("M64").Value = -4742
$ws.Range("N64").Value = -4322.6667

$ws.Range("H67").Value = 4292
$ws.Range("I67").Value = 4990
$ws.Range("J67").Value = 3826.6667
$ws.Range("K67").Value = 4990
$ws.Range("L67").Value = 3826.6667
$ws.Range("M67").Value = -4132
$ws.Range("N67").Value = -5542.6667

$ws.Range("H86").Value = 1722
$ws.Range("I86").Value = 1750.5454
$ws.Range("K86").Value = 1750.5454
$ws.Range("M86").Value = -627.5454

$ws.Range("H89").Value = 1722
$ws.Range("I89").Value = 1750.5454
$ws.Range("K89").Value = 8752.726999999999
$ws.Range("M89").Value = -3136.726999999999

$ws.Range("H113").Value = 4410.846
$ws.Range("I113").Value = 4049.375
$ws.Range("K113").Value = 4049.375
$ws.Range("M113").Value = -795.375

$ws.Range("H131").Value = 2968.75
$ws.Range("I131").Value = 2364.4119
$ws.Range("J131").Value = 4436.4287
$ws.Range("K131").Value = 7093.2357
$ws.Range("L131").Value = 13309.2861
$ws.Range("M131").Value = -2053.2357
$ws.Range("N131").Value = -23389.2861

$ws.Range("H137").Value = 5293.846
$ws.Range("I137").Value = 5748.5557
$ws.Range("J137").Value = 4270.75
$ws.Range("K137").Value = 17245.6671
$ws.Range("L137").Value = 12812.25
$ws.Range("M137").Value = -14695.6671
$ws.Range("N137").Value = -17912.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 2750
$ws.Range("J63").Value = 2750
$ws.Range("L63").Value = 2750
$ws.Range("N63").Value = -4122

$ws.Range("H66").Value = 2750
$ws.Range("J66").Value = 2750
$ws.Range("L66").Value = 13750
$ws.Range("N66").Value = -20614

$ws.Range("H122").Value = 2282.6943
$ws.Range("I122").Value = 1714.7587
$ws.Range("J122").Value = 4635.5713
$ws.Range("K122").Value = 5144.2761
$ws.Range("L122").Value = 13906.7139
$ws.Range("M122").Value = -2694.2761
$ws.Range("N122").Value = -18806.7139

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("N59").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 4646.8
$ws.Range("I62").Value = 2575
$ws.Range("J62").Value = 7754.5
$ws.Range("K62").Value = 2575
$ws.Range("L62").Value = 7754.5
$ws.Range("M62").Value = -1951
$ws.Range("N62").Value = -9002.5

$ws.Range("H65").Value = 4646.8
$ws.Range("I65").Value = 2575
$ws.Range("J65").Value = 7754.5
$ws.Range("K65").Value = 12875
$ws.Range("L65").Value = 38772.5
$ws.Range("M65").Value = -9755
$ws.Range("N65").Value = -45012.5

$ws.Range("H107").Value = 1117.075
$ws.Range("I107").Value = 989.069
$ws.Range("J107").Value = 1454.5454
$ws.Range("K107").Value = 989.069
$ws.Range("L107").Value = 1454.5454
$ws.Range("M107").Value = 930.931
$ws.Range("N107").Value = -5294.5454

$ws.Range("H122").Value = 2766.7778
$ws.Range("I122").Value = 2413.0476
$ws.Range("J122").Value = 4004.8333
$ws.Range("K122").Value = 7239.1428
$ws.Range("L122").Value = 12014.4999
$ws.Range("M122").Value = -4789.1428
$ws.Range("N122").Value = -16914.4999

$ws.Range("H134").Value = 9092888
$ws.Range("I134").Value = 13159498
$ws.Range("J134").Value = 2818.0588
$ws.Range("K134").Value = 39478494
$ws.Range("L134").Value = 8454.1764
$ws.Range("M134").Value = -39475959
$ws.Range("N134").Value = -13524.1764

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 372.5
$ws.Range("I14").Value = 372.5
$ws.Range("K14").Value = 1117.5
$ws.Range("M14").Value = -944.5

$ws.Range("H49").Value = 3650
$ws.Range("J49").Value = 3650
$ws.Range("L49").Value = 10950
$ws.Range("N49").Value = -11262

$ws.Range("H76").Value = 3666.6667
$ws.Range("I76").Value = 3000
$ws.Range("K76").Value = 9000
$ws.Range("M76").Value = -8617

$ws.Range("H79").Value = 3666.6667
$ws.Range("I79").Value = 3000
$ws.Range("K79").Value = 9000
$ws.Range("M79").Value = -7674

$ws.Range("H82").Value = 3182.6
$ws.Range("I82").Value = 1956.5
$ws.Range("K82").Value = 5869.5
$ws.Range("M82").Value = -5463.5

$ws.Range("H85").Value = 3182.6
$ws.Range("I85").Value = 1956.5
$ws.Range("K85").Value = 5869.5
$ws.Range("M85").Value = -4465.5

$ws.Range("H107").Value = 1588.2858
$ws.Range("I107").Value = 366.4
$ws.Range("J107").Value = 2267.111
$ws.Range("K107").Value = 1099.2
$ws.Range("L107").Value = 6801.333
$ws.Range("M107").Value = 820.8000000000002
$ws.Range("N107").Value = -10641.333

$ws.Range("H119").Value = 3694.5386
$ws.Range("I119").Value = 2014.5
$ws.Range("K119").Value = 6043.5
$ws.Range("M119").Value = -1205.5

$ws.Range("H122").Value = 1648.625
$ws.Range("I122").Value = 550.25
$ws.Range("J122").Value = 2747
$ws.Range("K122").Value = 4952.25
$ws.Range("L122").Value = 24723
$ws.Range("M122").Value = -2502.25
$ws.Range("N122").Value = -29623

$ws.Range("H131").Value = 1508.4783
$ws.Range("I131").Value = 3681
$ws.Range("J131").Value = 1118.5385
$ws.Range("K131").Value = 11043
$ws.Range("L131").Value = 3355.6155
$ws.Range("M131").Value = -6003
$ws.Range("N131").Value = -13435.6155

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 1264.3334
$ws.Range("I55").Value = 235.25
$ws.Range("K55").Value = 235.25
$ws.Range("M55").Value = -62.25

$ws.Range("H68").Value = 1612
$ws.Range("I68").Value = 1170.5264
$ws.Range("J68").Value = 10000
$ws.Range("K68").Value = 1170.5264
$ws.Range("L68").Value = 10000
$ws.Range("M68").Value = -421.5264
$ws.Range("N68").Value = -11498

$ws.Range("H71").Value = 1612
$ws.Range("I71").Value = 1170.5264
$ws.Range("J71").Value = 10000
$ws.Range("K71").Value = 5852.632
$ws.Range("L71").Value = 50000
$ws.Range("M71").Value = -2108.632
$ws.Range("N71").Value = -57488

$ws.Range("H139").Value = 24172.69
$ws.Range("J139").Value = 24172.69
$ws.Range("L139").Value = 24172.69
$ws.Range("N139").Value = -34452.69

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H55").Value = 5499.5
$ws.Range("I55").Value = 3000
$ws.Range("J55").Value = 7999
$ws.Range("K55").Value = 3000
$ws.Range("L55").Value = 7999
$ws.Range("M55").Value = -2723
$ws.Range("N55").Value = -8553

$ws.Range("H62").Value = 3500
$ws.Range("I62").Value = 3000
$ws.Range("J62").Value = 4000
$ws.Range("K62").Value = 3000
$ws.Range("L62").Value = 4000
$ws.Range("M62").Value = -2376
$ws.Range("N62").Value = -5248

$ws.Range("H65").Value = 3500
$ws.Range("I65").Value = 3000
$ws.Range("J65").Value = 4000
$ws.Range("K65").Value = 15000
$ws.Range("L65").Value = 20000
$ws.Range("M65").Value = -11880
$ws.Range("N65").Value = -26240

$ws.Range("H113").Value = 1427.125
$ws.Range("I113").Value = 183.4
$ws.Range("J113").Value = 3500
$ws.Range("K113").Value = 550.2
$ws.Range("L113").Value = 10500
$ws.Range("M113").Value = 1619.8
$ws.Range("N113").Value = -14840

$ws.Range("H136").Value = 2193.152
$ws.Range("I136").Value = 1842.8276
$ws.Range("J136").Value = 2790.7646
$ws.Range("K136").Value = 5528.4828
$ws.Range("L136").Value = 8372.293799999999
$ws.Range("M136").Value = -2978.4828
$ws.Range("N136").Value = -13472.2938
